$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update region names in column B (지역) so map view labels show the
# full administrative name instead of the short city/county name:
#   영천 -> 영천시, 임실 -> 임실군, 고흥 -> 고흥군, 문경 -> 문경시
$ws.Range("B2:B6").Value = "영천시"
$ws.Range("B7").Value = "임실군"
$ws.Range("B8").Value = "고흥군"
$ws.Range("B9").Value = "문경시"
